# New crime data collected — weekly refresh of the 45th Precinct CompStat
# report: bump the report Volume/Number and the covered week's date range,
# then refresh the week/28-day/YTD/2-yr/13-yr/30-yr crime-complaint figures
# for every offense row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: some cells in the stats grid hold the literal text "0" or
# "***.*" as placeholders (shown when a row has no prior-year base to
# compute a % change against). Assigning a bare numeric-looking string
# to .Value lets the engine coerce it back to a real number, so we set
# the value with a leading apostrophe (forces text) and then repaint
# the cell's format from a donor cell that already carries the correct
# "text" style, restoring the exact original look.
# ---------------------------------------------------------------------
function Set-PlaceholderText($ref, $text, $donorRef) {
    $ws.Range($ref).Value = "'" + $text
    $ws.Range($donorRef).Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

# --- Header: Volume/Number and the reporting week's date range --------
$ws.Range("A8").Value = "Volume 30   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/30/2023  Through  11/5/2023"

# --- Murder (row 14) ----------------------------------------------------
Set-PlaceholderText "F14" "0" "C14"

# --- Rape (row 15) -------------------------------------------------------
$ws.Range("C15").Value = 1
$ws.Range("I15").Value = 17
$ws.Range("K15").Value = 41.666666666666
$ws.Range("L15").Value = 41.666666666666
$ws.Range("M15").Value = -5.555555555555
$ws.Range("N15").Value = -15

# --- Robbery (row 16) -----------------------------------------------------
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = -31.818181818181
$ws.Range("I16").Value = 175
$ws.Range("J16").Value = 177
$ws.Range("K16").Value = -1.129943502824
$ws.Range("L16").Value = 19.047619047619
$ws.Range("M16").Value = -10.25641025641
$ws.Range("N16").Value = -55.919395465995

# --- Fel. Assault (row 17) -------------------------------------------------
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -14.285714285714
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -11.764705882352
$ws.Range("I17").Value = 252
$ws.Range("J17").Value = 222
$ws.Range("K17").Value = 13.513513513513
$ws.Range("L17").Value = 29.230769230769
$ws.Range("M17").Value = 58.490566037735
$ws.Range("N17").Value = 26.633165829145

# --- Burglary (row 18) ------------------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -10
$ws.Range("I18").Value = 119
$ws.Range("J18").Value = 89
$ws.Range("K18").Value = 33.707865168539
$ws.Range("L18").Value = 25.263157894736
$ws.Range("M18").Value = -54.406130268199
$ws.Range("N18").Value = -83.356643356643

# --- Gr. Larceny (row 19) ---------------------------------------------------
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 175
$ws.Range("F19").Value = 62
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = 16.981132075471
$ws.Range("I19").Value = 564
$ws.Range("J19").Value = 481
$ws.Range("K19").Value = 17.255717255717
$ws.Range("L19").Value = 41
$ws.Range("M19").Value = 37.226277372262
$ws.Range("N19").Value = 51.612903225806

# --- G.L.A. (row 20) ---------------------------------------------------------
$ws.Range("C20").Value = 6
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 35
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = 20.689655172413
$ws.Range("I20").Value = 445
$ws.Range("J20").Value = 258
$ws.Range("K20").Value = 72.480620155038
$ws.Range("L20").Value = 72.480620155038
$ws.Range("M20").Value = 150
$ws.Range("N20").Value = -75.656455142231

# --- TOTAL (row 21) -----------------------------------------------------------
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = 37.037037037037
$ws.Range("F21").Value = 140
$ws.Range("G21").Value = 131
$ws.Range("H21").Value = 6.870229007633
$ws.Range("I21").Value = 1576
$ws.Range("J21").Value = 1243
$ws.Range("K21").Value = 26.790024135156
$ws.Range("L21").Value = 41.726618705036
$ws.Range("M21").Value = 28.968903436988
$ws.Range("N21").Value = -55.530474040632

# --- Transit (row 22) -----------------------------------------------------------
$ws.Range("C22").Value = 3
$ws.Range("F22").Value = 5
$ws.Range("I22").Value = 16
$ws.Range("K22").Value = 23.076923076923
$ws.Range("L22").Value = 23.076923076923
$ws.Range("M22").Value = 33.333333333333

# --- Housing (row 23) -----------------------------------------------------------
Set-PlaceholderText "C23" "0" "D22"
Set-PlaceholderText "D23" "0" "D22"
Set-PlaceholderText "E23" "***.*" "E22"
$ws.Range("F23").Value = 4
$ws.Range("H23").Value = 100
$ws.Range("L23").Value = 29.268292682926
$ws.Range("M23").Value = 15.217391304347

# --- Petit Larceny (row 24) -------------------------------------------------------
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = -5.882352941176
$ws.Range("F24").Value = 92
$ws.Range("G24").Value = 91
$ws.Range("H24").Value = 1.098901098901
$ws.Range("I24").Value = 1181
$ws.Range("J24").Value = 985
$ws.Range("K24").Value = 19.89847715736
$ws.Range("L24").Value = 43.325242718446
$ws.Range("M24").Value = -6.78768745067

# --- Misd. Assault (row 25) ---------------------------------------------------------
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = 9.756097560975
$ws.Range("I25").Value = 437
$ws.Range("J25").Value = 393
$ws.Range("K25").Value = 11.19592875318
$ws.Range("L25").Value = 16.223404255319
$ws.Range("M25").Value = 18.108108108108

# --- UCR Rape* (row 26) ---------------------------------------------------------------
$ws.Range("C26").Value = 1
$ws.Range("I26").Value = 30
$ws.Range("K26").Value = 15.384615384615
$ws.Range("L26").Value = 87.5

# --- Other Sex Crimes (row 27) ---------------------------------------------------------
# C27 used to be a blank "0" placeholder (text style); this week it has a
# real count, so give it the donor-cell numeric format before writing the
# number (mirrors Set-PlaceholderText's format-then-paste trick, just
# landing on the numeric style instead of the text one).
$ws.Range("C27").Value = 3
$ws.Range("C16").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
Set-PlaceholderText "D27" "0" "D22"
Set-PlaceholderText "E27" "***.*" "E22"
$ws.Range("F27").Value = 6
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 39
$ws.Range("K27").Value = -13.333333333333
$ws.Range("L27").Value = 5.405405405405

# --- Shooting Vic. (row 28) ---------------------------------------------------------------
$ws.Range("L28").Value = -18.181818181818

# --- Shooting Inc. (row 29) ---------------------------------------------------------------
$ws.Range("L29").Value = -22.222222222222

# --- Hate Crimes (row 30) ---------------------------------------------------------------
Set-PlaceholderText "D30" "0" "D22"
Set-PlaceholderText "E30" "***.*" "E22"
